# Weekly fruit/vegetable price update: insert one new observation row
# (a new "Primera" quality record dated 2021-12-24) right before the
# existing row 41, shifting every following record down by one row
# (old row 41 becomes row 42, ... old row 85 becomes row 86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41; pushes rows 41..85 down to 42..86
# and copies the formatting of the row above (-4121 = xlShiftDown).
$ws.Rows(41).Insert(-4121)

# Populate the newly inserted row 41 with the new observation.
$ws.Cells.Item(41, 1).Value  = 1
$ws.Cells.Item(41, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(41, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(41, 4).Value  = 44554
$ws.Cells.Item(41, 5).Value  = 15
$ws.Cells.Item(41, 6).Value  = 100112036
$ws.Cells.Item(41, 7).Value  = "Caigua"
$ws.Cells.Item(41, 8).Value  = "Sin especificar"
$ws.Cells.Item(41, 9).Value  = "Primera"
$ws.Cells.Item(41, 10).Value = 160
$ws.Cells.Item(41, 11).Value = 6000
$ws.Cells.Item(41, 12).Value = 7000
$ws.Cells.Item(41, 13).Value = 6500
$ws.Cells.Item(41, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(41, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(41, 16).Value = 325
$ws.Cells.Item(41, 17).Value = 20
$ws.Cells.Item(41, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format as the rest
# of column D.
$ws.Cells.Item(41, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
